# Homehelp status dates refactoring
#
# The "Ujra nyitott" (reopened) status/date row was removed from the
# catering_stats report. It used to sit at row 5 (between "Teljes havi
# letszam" and "Aktiv ellatottak szama"); deleting it shifts every
# subsequent row up by one, which also renumbers/repacks the merged
# cell ranges and drops the now-unused shared strings automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("5").Delete()
